$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to remain a text value (matches the source data's
    # inline-string storage) instead of letting Excel auto-coerce
    # numeric-looking strings ("0.999", "132.43", ...) into numbers.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.947.72"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.352.33"
$ws.Range("E3").Value = "  +0.02%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "550.53"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6 - Solana
Set-TextValue "D6" "132.43"
$ws.Range("E6").Value = "  -1.93%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.71%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +3.41%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +4.30%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.22%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.353"
$ws.Range("E12").Value = "  -1.61%  "

# Row 13 - Avalanche
Set-TextValue "D13" "24.01"
$ws.Range("E13").Value = "  +1.79%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.769.79"
$ws.Range("E14").Value = "  +0.01%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "57.832.30"
$ws.Range("E15").Value = "  +0.14%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.84%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.330.04"
$ws.Range("E17").Value = "  -1.10%  "

# Row 18 - Chainlink
Set-TextValue "D18" "10.98"
$ws.Range("E18").Value = "  +3.09%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.41%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "329.95"
$ws.Range("E20").Value = "  -1.47%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.87"
$ws.Range("E21").Value = "  +2.23%  "

# Row 22 - Dai
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23 - Litecoin
Set-TextValue "D23" "63.84"
$ws.Range("E23").Value = "  +2.71%  "

# Row 24 - Kaspa
$ws.Range("E24").Value = "  -0.17%  "

# Row 25 - Binance-PegBSC-USD
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  -0.11%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "8.26"
$ws.Range("E26").Value = "  -2.80%  "

# Row 27 - Fetch.AI
Set-TextValue "D27" "1.32"
$ws.Range("E27").Value = "  -5.48%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "1.77"
$ws.Range("E28").Value = "  -0.49%  "

# Row 29 - Monero
Set-TextValue "D29" "171.05"
$ws.Range("E29").Value = "  +0.50%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31 - Aptos
Set-TextValue "D31" "6.13"
$ws.Range("E31").Value = "  -0.45%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "18.34"
$ws.Range("E32").Value = "  -1.02%  "

# Row 33 - SuiNetwork
$ws.Range("E33").Value = "  -2.67%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "4.14"
$ws.Range("E36").Value = "  -0.99%  "

# Row 37 - was ImmutableX, now PolygonEcosystemToken
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D37" "0.437"
$ws.Range("E37").Value = "  +15.96%  "

# Row 38 - was PolygonEcosystemToken, now ImmutableX
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "1.24"
$ws.Range("E38").Value = "  -1.36%  "

# Row 39 - OKB
Set-TextValue "D39" "40.36"
$ws.Range("E39").Value = "  +3.14%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.60"
$ws.Range("E40").Value = "  -1.22%  "

# Row 41 - Aave
Set-TextValue "D41" "142.08"
$ws.Range("E41").Value = "  -3.92%  "

# Row 42 - Filecoin
Set-TextValue "D42" "3.65"
$ws.Range("E42").Value = "  +0.36%  "

# Row 43 - Bittensor
Set-TextValue "D43" "288.63"
$ws.Range("E43").Value = "  +1.11%  "

# Row 44 - Polygon
Set-TextValue "D44" "0.426"
$ws.Range("E44").Value = "  +10.74%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +0.97%  "

# Row 46 - Hedera
Set-TextValue "D46" "0.0514"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47 - Mantle
Set-TextValue "D47" "0.566"
$ws.Range("E47").Value = "  +0.81%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "18.75"
$ws.Range("E48").Value = "  -2.45%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0222"
$ws.Range("E49").Value = "  +1.68%  "

# Row 50 - WhiteBITCoin
$ws.Range("E50").Value = "  +0.01%  "
